# Apply updated crypto price/volume data to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 4).Value = "48.024.02"
$ws.Cells.Item(2, 5).Value = "  +2.38%  "

$ws.Cells.Item(3, 4).Value = "2.515.75"
$ws.Cells.Item(3, 5).Value = "  +1.72%  "

$ws.Cells.Item(4, 5).Value = "  -0.05%  "

$ws.Cells.Item(5, 4).Value = "'321.78"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +0.87%  "

$ws.Cells.Item(6, 4).Value = "'108.81"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  +1.38%  "

$ws.Cells.Item(7, 4).Value = "'0.531"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  +2.35%  "

$ws.Cells.Item(8, 5).Value = "  +0.02%  "

$ws.Cells.Item(9, 4).Value = "'0.552"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  +4.31%  "

$ws.Cells.Item(10, 4).Value = "'40.35"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  +5.49%  "

$ws.Cells.Item(11, 5).Value = "  +14.10%  "

$ws.Cells.Item(12, 5).Value = "  +2.41%  "

$ws.Cells.Item(13, 5).Value = "  +1.35%  "

$ws.Cells.Item(14, 4).Value = "'7.25"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  +2.53%  "

$ws.Cells.Item(15, 4).Value = "2.911.52"
$ws.Cells.Item(15, 5).Value = "  +1.54%  "

$ws.Cells.Item(16, 4).Value = "2.517.83"
$ws.Cells.Item(16, 5).Value = "  +1.22%  "

$ws.Cells.Item(17, 4).Value = "'0.854"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  +1.93%  "

$ws.Cells.Item(18, 4).Value = "47.876.02"
$ws.Cells.Item(18, 5).Value = "  +2.13%  "

$ws.Cells.Item(19, 4).Value = "'13.25"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  +5.07%  "

$ws.Cells.Item(20, 4).Value = "'6.61"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  +0.91%  "

$ws.Cells.Item(21, 4).Value = "0.0₃0944"
$ws.Cells.Item(21, 5).Value = "  +2.13%  "

$ws.Cells.Item(22, 4).Value = "'2.70"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  -1.26%  "

$ws.Cells.Item(23, 4).Value = "'71.88"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  +2.32%  "

$ws.Cells.Item(24, 4).Value = "'264.49"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  +8.54%  "

$ws.Cells.Item(25, 4).Value = "'2.56"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  +1.13%  "

$ws.Cells.Item(26, 5).Value = "  -0.22%  "

$ws.Cells.Item(27, 4).Value = "'26.04"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  +2.42%  "

$ws.Cells.Item(28, 5).Value = "  +1.58%  "

$ws.Cells.Item(29, 2).Value = "Toncoin"
$ws.Cells.Item(29, 3).Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Cells.Item(29, 4).Value = "'2.20"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  +0.90%  "

$ws.Cells.Item(30, 2).Value = "Kaspa"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(30, 4).Value = "'0.143"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  +2.47%  "

$ws.Cells.Item(31, 4).Value = "'35.88"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  +4.33%  "

$ws.Cells.Item(32, 4).Value = "'49.71"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  +0.47%  "

$ws.Cells.Item(33, 5).Value = "  -0.01%  "

$ws.Cells.Item(34, 4).Value = "'5.39"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  +1.99%  "

$ws.Cells.Item(35, 5).Value = "  -0.13%  "

$ws.Cells.Item(36, 4).Value = "'0.0788"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  +1.75%  "

$ws.Cells.Item(37, 5).Value = "  +2.03%  "

$ws.Cells.Item(38, 4).Value = "'4.69"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  +2.22%  "

$ws.Cells.Item(39, 5).Value = "  +2.76%  "

$ws.Cells.Item(40, 5).Value = "  +0.90%  "

$ws.Cells.Item(41, 4).Value = "'22.08"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  +4.33%  "

$ws.Cells.Item(42, 2).Value = "Monero"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(42, 4).Value = "'119.65"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  +0.82%  "

$ws.Cells.Item(43, 2).Value = "WEMIXToken"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Cells.Item(43, 4).Value = "'2.20"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  -0.57%  "

$ws.Cells.Item(44, 5).Value = "  +2.83%  "

$ws.Cells.Item(45, 4).Value = "2.010.79"
$ws.Cells.Item(45, 5).Value = "  +2.20%  "

$ws.Cells.Item(46, 4).Value = "'3.15"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  +6.08%  "

$ws.Cells.Item(47, 4).Value = "'1.90"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  +9.38%  "

$ws.Cells.Item(48, 5).Value = "  +2.86%  "

$ws.Cells.Item(49, 4).Value = "'9.09"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  +1.14%  "

$ws.Cells.Item(50, 4).Value = "'5.19"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  +2.05%  "

$ws.Cells.Item(51, 4).Value = "'78.59"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  +3.19%  "
